$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1240.125
$ws.Range("I15").Value = 1240.125
$ws.Range("K15").Value = 3720.375
$ws.Range("M15").Value = -3551.375
$ws.Range("H29").Value = 1124.5
$ws.Range("I29").Value = 999.3333
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 2997.9999
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = -2716.9999
$ws.Range("N29").Value = -5062
$ws.Range("H32").Value = 8276.6
$ws.Range("J32").Value = 11997.833
$ws.Range("L32").Value = 11997.833
$ws.Range("N32").Value = -12649.833
$ws.Range("H116").Value = 3499.2
$ws.Range("I116").Value = 3499.25
$ws.Range("J116").Value = 3499
$ws.Range("K116").Value = 3499.25
$ws.Range("L116").Value = 3499
$ws.Range("M116").Value = -57.25
$ws.Range("N116").Value = -10383
$ws.Range("H137").Value = 6549.1
$ws.Range("I137").Value = 6540.9287
$ws.Range("J137").Value = 6568.1665
$ws.Range("K137").Value = 19622.7861
$ws.Range("L137").Value = 19704.4995
$ws.Range("M137").Value = -17072.7861
$ws.Range("N137").Value = -24804.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1640.3684
$ws.Range("I32").Value = 1640.3684
$ws.Range("K32").Value = 1640.3684
$ws.Range("M32").Value = -1353.3684
$ws.Range("H102").Value = 536.93335
$ws.Range("I102").Value = 519.53845
$ws.Range("K102").Value = 519.53845
$ws.Range("M102").Value = 1102.46155
$ws.Range("H122").Value = 2222
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2222
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6666
$ws.Range("N122").Value = -11566
$ws.Range("H132").Value = 4095.4546
$ws.Range("I132").Value = 4399.6113
$ws.Range("K132").Value = 13198.8339
$ws.Range("M132").Value = -10668.8339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7298.6
$ws.Range("I86").Value = 3426.8572
$ws.Range("K86").Value = 3426.8572
$ws.Range("M86").Value = -2303.8572
$ws.Range("H89").Value = 7298.6
$ws.Range("I89").Value = 3426.8572
$ws.Range("K89").Value = 17134.286
$ws.Range("M89").Value = -11518.286
$ws.Range("H105").Value = 6058.1113
$ws.Range("I105").Value = 5769.7334
$ws.Range("K105").Value = 5769.7334
$ws.Range("M105").Value = -4022.7334
$ws.Range("H134").Value = 8121.6665
$ws.Range("I134").Value = 7907.8335
$ws.Range("K134").Value = 23723.5005
$ws.Range("M134").Value = -21188.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1937.75
$ws.Range("I6").Value = 1937.75
$ws.Range("K6").Value = 1937.75
$ws.Range("M6").Value = -1824.75
$ws.Range("H16").Value = 2030
$ws.Range("J16").Value = 3087.5
$ws.Range("L16").Value = 3087.5
$ws.Range("N16").Value = -3661.5
$ws.Range("H31").Value = 2991
$ws.Range("I31").Value = 2991
$ws.Range("K31").Value = 2991
$ws.Range("M31").Value = -2696
$ws.Range("H34").Value = 2991
$ws.Range("I34").Value = 2991
$ws.Range("K34").Value = 2991
$ws.Range("M34").Value = -2789
$ws.Range("H107").Value = 850.8
$ws.Range("I107").Value = 752
$ws.Range("K107").Value = 752
$ws.Range("M107").Value = 1168
$ws.Range("H113").Value = 2030
$ws.Range("J113").Value = 3087.5
$ws.Range("L113").Value = 3087.5
$ws.Range("N113").Value = -7427.5
$ws.Range("H141").Value = 169998.5
$ws.Range("J141").Value = 169998.5
$ws.Range("L141").Value = 169998.5
$ws.Range("N141").Value = -180358.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 82.875
$ws.Range("I6").Value = 59
$ws.Range("K6").Value = 177
$ws.Range("M6").Value = -64
$ws.Range("H18").Value = 320
$ws.Range("I18").Value = 266.25
$ws.Range("K18").Value = 798.75
$ws.Range("M18").Value = -629.75
$ws.Range("H34").Value = 77730.14
$ws.Range("J34").Value = 95990.17999999999
$ws.Range("L34").Value = 287970.54
$ws.Range("N34").Value = -288138.54
$ws.Range("H39").Value = 7187.5
$ws.Range("J39").Value = 7187.5
$ws.Range("L39").Value = 21562.5
$ws.Range("N39").Value = -22150.5
$ws.Range("H55").Value = 16873.25
$ws.Range("J55").Value = 16873.25
$ws.Range("L55").Value = 50619.75
$ws.Range("N55").Value = -50973.75
$ws.Range("H121").Value = 111112090
$ws.Range("J121").Value = 166667820
$ws.Range("L121").Value = 500003460
$ws.Range("N121").Value = -500006080
$ws.Range("H140").Value = 627181.25
$ws.Range("I140").Value = 627181.25
$ws.Range("K140").Value = 1881543.75
$ws.Range("M140").Value = -1876363.75
$ws.Range("H141").Value = 7268.909
$ws.Range("I141").Value = 7576.6665
$ws.Range("J141").Value = 6899.6
$ws.Range("K141").Value = 22729.9995
$ws.Range("L141").Value = 20698.8
$ws.Range("M141").Value = -17549.9995
$ws.Range("N141").Value = -31058.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2898.25
$ws.Range("I80").Value = 2898.6667
$ws.Range("J80").Value = 2897
$ws.Range("K80").Value = 2898.6667
$ws.Range("L80").Value = 2897
$ws.Range("M80").Value = -1900.6667
$ws.Range("N80").Value = -4893
$ws.Range("H83").Value = 2898.25
$ws.Range("I83").Value = 2898.6667
$ws.Range("J83").Value = 2897
$ws.Range("K83").Value = 14493.3335
$ws.Range("L83").Value = 14485
$ws.Range("M83").Value = -9501.333500000001
$ws.Range("N83").Value = -24469
$ws.Range("H95").Value = 28333.334
$ws.Range("J95").Value = 28333.334
$ws.Range("L95").Value = 28333.334
$ws.Range("N95").Value = -33825.334
$ws.Range("H122").Value = 3889
$ws.Range("I122").Value = 3889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11667
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9217
$ws.Range("H132").Value = 2156.3125
$ws.Range("I132").Value = 2170.923
$ws.Range("K132").Value = 6512.768999999999
$ws.Range("M132").Value = -3982.768999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1588.6
$ws.Range("I46").Value = 947.5714
$ws.Range("K46").Value = 947.5714
$ws.Range("M46").Value = -759.5714
$ws.Range("H82").Value = 1683.9286
$ws.Range("I82").Value = 1023.25
$ws.Range("J82").Value = 1948.2
$ws.Range("K82").Value = 1023.25
$ws.Range("L82").Value = 1948.2
$ws.Range("M82").Value = -662.25
$ws.Range("N82").Value = -2670.2
$ws.Range("H85").Value = 1683.9286
$ws.Range("I85").Value = 1023.25
$ws.Range("J85").Value = 1948.2
$ws.Range("K85").Value = 1023.25
$ws.Range("L85").Value = 1948.2
$ws.Range("M85").Value = 224.75
$ws.Range("N85").Value = -4444.2
$ws.Range("H93").Value = 1364.9
$ws.Range("I93").Value = 1468.625
$ws.Range("K93").Value = 1468.625
$ws.Range("M93").Value = -220.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2425
$ws.Range("I100").Value = 468
$ws.Range("K100").Value = 936
$ws.Range("M100").Value = -395
$ws.Range("H113").Value = 823.0769
$ws.Range("I113").Value = 518.1818
$ws.Range("K113").Value = 1554.5454
$ws.Range("M113").Value = 615.4546
$ws.Range("H136").Value = 6901.32
$ws.Range("I136").Value = 6901.32
$ws.Range("K136").Value = 20703.96
$ws.Range("M136").Value = -18153.96

# Cell removals (entirely clear the cell, matching the authoritative diff)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N122").ClearContents()
